$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 55, shifting existing rows (55-76) down to (56-77)
$ws.Rows(55).Insert()

# Helper-free approach: force text storage for numeric-looking strings by
# temporarily applying a text number format, then clearing the format so the
# cell ends up with the default (unstyled) appearance, matching the rest of
# the sheet's data rows.

$ws.Range("A55").NumberFormat = "@"
$ws.Range("A55").Value = "6002"
$ws.Range("A55").ClearFormats()

$ws.Range("B55").NumberFormat = "@"
$ws.Range("B55").Value = "8/12/2025"
$ws.Range("B55").ClearFormats()

$ws.Range("C55").Value = "LA PLATA AV. 832"

$ws.Range("D55").NumberFormat = "@"
$ws.Range("D55").Value = "5"
$ws.Range("D55").ClearFormats()

$ws.Range("E55").NumberFormat = "@"
$ws.Range("E55").Value = "808918694"
$ws.Range("E55").ClearFormats()

$ws.Range("F55").Value = "NEW"
$ws.Range("G55").Value = "Pendiente"
$ws.Range("H55").Value = "Picada"
$ws.Range("I55").Value = 1
$ws.Range("J55").Value = "Cambio"
$ws.Range("K55").Value = "Sin equipos"
$ws.Range("L55").Value = "Pasante"
$ws.Range("M55").Value = -58.426947
$ws.Range("N55").Value = -34.625698
$ws.Range("O55").Value = "Boedo"
$ws.Range("P55").Value = "Capital Sur"
